$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student name in row 3: "Filipe Palma Abreu" -> "Lucas Borges Jagersbacher"
$ws.Range("A3").Value = "Lucas Borges Jagersbacher"

# Update the grades for row 2 (Eduardo Lago Nunes): set R3 (D2) grade
$ws.Range("D2").Value = 0

# Update the grades for row 3 (now Lucas Borges Jagersbacher): R2 (C3) and R3 (D3)
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2

# Update the grades for row 4 (Rosa Maria Araujo Calazans): set R3 (D4) grade
$ws.Range("D4").Value = 0

# Update the selection to A4, as in the saved workbook view
$ws.Range("A4").Select()
